# Populate the previously-empty Circulation / ILL Loans / ILL Borrows
# figures (columns B:D) for rows 3-59 of Sheet1, matching the
# "January 1-August, 31, 2022" statistics update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 2).Value = 50519
$ws.Cells.Item(3, 3).Value = 7195
$ws.Cells.Item(3, 4).Value = 8169
$ws.Cells.Item(4, 2).Value = 23927
$ws.Cells.Item(4, 3).Value = 2610
$ws.Cells.Item(4, 4).Value = 2886
$ws.Cells.Item(5, 2).Value = 84446
$ws.Cells.Item(5, 3).Value = 5819
$ws.Cells.Item(5, 4).Value = 7580
$ws.Cells.Item(6, 2).Value = 1832
$ws.Cells.Item(6, 3).Value = 520
$ws.Cells.Item(6, 4).Value = 73
$ws.Cells.Item(7, 2).Value = 49022
$ws.Cells.Item(7, 3).Value = 8499
$ws.Cells.Item(7, 4).Value = 7414
$ws.Cells.Item(8, 2).Value = 5735
$ws.Cells.Item(8, 3).Value = 1237
$ws.Cells.Item(8, 4).Value = 1201
$ws.Cells.Item(9, 2).Value = 6474
$ws.Cells.Item(9, 3).Value = 1103
$ws.Cells.Item(9, 4).Value = 845
$ws.Cells.Item(10, 2).Value = 2488
$ws.Cells.Item(10, 3).Value = 355
$ws.Cells.Item(10, 4).Value = 92
$ws.Cells.Item(11, 2).Value = 408
$ws.Cells.Item(11, 3).Value = 263
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(13, 2).Value = 1147
$ws.Cells.Item(13, 3).Value = 335
$ws.Cells.Item(13, 4).Value = 252
$ws.Cells.Item(14, 2).Value = 3347
$ws.Cells.Item(14, 3).Value = 1534
$ws.Cells.Item(14, 4).Value = 998
$ws.Cells.Item(15, 2).Value = 6776
$ws.Cells.Item(15, 3).Value = 2297
$ws.Cells.Item(15, 4).Value = 1078
$ws.Cells.Item(16, 2).Value = 3860
$ws.Cells.Item(16, 3).Value = 1634
$ws.Cells.Item(16, 4).Value = 676
$ws.Cells.Item(17, 2).Value = 2214
$ws.Cells.Item(17, 3).Value = 724
$ws.Cells.Item(17, 4).Value = 170
$ws.Cells.Item(18, 2).Value = 19670
$ws.Cells.Item(18, 3).Value = 2998
$ws.Cells.Item(18, 4).Value = 3508
$ws.Cells.Item(19, 2).Value = 3200
$ws.Cells.Item(19, 3).Value = 886
$ws.Cells.Item(19, 4).Value = 566
$ws.Cells.Item(20, 2).Value = 24006
$ws.Cells.Item(20, 3).Value = 2556
$ws.Cells.Item(20, 4).Value = 4099
$ws.Cells.Item(21, 2).Value = 373
$ws.Cells.Item(21, 3).Value = 385
$ws.Cells.Item(21, 4).Value = 14
$ws.Cells.Item(22, 2).Value = 20157
$ws.Cells.Item(22, 3).Value = 2388
$ws.Cells.Item(22, 4).Value = 3111
$ws.Cells.Item(23, 2).Value = 1356
$ws.Cells.Item(23, 3).Value = 541
$ws.Cells.Item(23, 4).Value = 206
$ws.Cells.Item(24, 2).Value = 20756
$ws.Cells.Item(24, 3).Value = 3532
$ws.Cells.Item(24, 4).Value = 3543
$ws.Cells.Item(25, 2).Value = 77931
$ws.Cells.Item(25, 3).Value = 8919
$ws.Cells.Item(25, 4).Value = 8646
$ws.Cells.Item(26, 2).Value = 5827
$ws.Cells.Item(26, 3).Value = 1929
$ws.Cells.Item(26, 4).Value = 792
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(28, 2).Value = 6429
$ws.Cells.Item(28, 3).Value = 1244
$ws.Cells.Item(28, 4).Value = 1255
$ws.Cells.Item(29, 2).Value = 2714
$ws.Cells.Item(29, 3).Value = 380
$ws.Cells.Item(29, 4).Value = 546
$ws.Cells.Item(30, 2).Value = 17860
$ws.Cells.Item(30, 3).Value = 3199
$ws.Cells.Item(30, 4).Value = 3026
$ws.Cells.Item(31, 2).Value = 694
$ws.Cells.Item(31, 3).Value = 260
$ws.Cells.Item(31, 4).Value = 243
$ws.Cells.Item(32, 2).Value = 2874
$ws.Cells.Item(32, 3).Value = 1811
$ws.Cells.Item(32, 4).Value = 325
$ws.Cells.Item(33, 2).Value = 17047
$ws.Cells.Item(33, 3).Value = 3387
$ws.Cells.Item(33, 4).Value = 3061
$ws.Cells.Item(34, 2).Value = 12877
$ws.Cells.Item(34, 3).Value = 2926
$ws.Cells.Item(34, 4).Value = 3100
$ws.Cells.Item(35, 2).Value = 7261
$ws.Cells.Item(35, 3).Value = 861
$ws.Cells.Item(35, 4).Value = 1518
$ws.Cells.Item(36, 2).Value = 60202
$ws.Cells.Item(36, 3).Value = 6439
$ws.Cells.Item(36, 4).Value = 6187
$ws.Cells.Item(37, 2).Value = 9693
$ws.Cells.Item(37, 3).Value = 2928
$ws.Cells.Item(37, 4).Value = 1527
$ws.Cells.Item(38, 2).Value = 22790
$ws.Cells.Item(38, 3).Value = 2161
$ws.Cells.Item(38, 4).Value = 3240
$ws.Cells.Item(39, 2).Value = 1053
$ws.Cells.Item(39, 3).Value = 852
$ws.Cells.Item(39, 4).Value = 178
$ws.Cells.Item(40, 2).Value = 1810
$ws.Cells.Item(40, 3).Value = 299
$ws.Cells.Item(40, 4).Value = 802
$ws.Cells.Item(41, 2).Value = 2237
$ws.Cells.Item(41, 3).Value = 272
$ws.Cells.Item(41, 4).Value = 130
$ws.Cells.Item(42, 2).Value = 9137
$ws.Cells.Item(42, 3).Value = 255
$ws.Cells.Item(42, 4).Value = 198
$ws.Cells.Item(43, 2).Value = 298
$ws.Cells.Item(43, 3).Value = 113
$ws.Cells.Item(43, 4).Value = 78
$ws.Cells.Item(44, 2).Value = 761
$ws.Cells.Item(44, 3).Value = 53
$ws.Cells.Item(44, 4).Value = 22
$ws.Cells.Item(45, 2).Value = 1980
$ws.Cells.Item(45, 3).Value = 159
$ws.Cells.Item(45, 4).Value = 73
$ws.Cells.Item(46, 2).Value = 3711
$ws.Cells.Item(46, 3).Value = 1053
$ws.Cells.Item(46, 4).Value = 543
$ws.Cells.Item(47, 2).Value = 13809
$ws.Cells.Item(47, 3).Value = 3568
$ws.Cells.Item(47, 4).Value = 2682
$ws.Cells.Item(48, 2).Value = 34394
$ws.Cells.Item(48, 3).Value = 3495
$ws.Cells.Item(48, 4).Value = 5153
$ws.Cells.Item(49, 2).Value = 17334
$ws.Cells.Item(49, 3).Value = 3604
$ws.Cells.Item(49, 4).Value = 1343
$ws.Cells.Item(50, 2).Value = 12282
$ws.Cells.Item(50, 3).Value = 1201
$ws.Cells.Item(50, 4).Value = 2050
$ws.Cells.Item(51, 2).Value = 31505
$ws.Cells.Item(51, 3).Value = 3239
$ws.Cells.Item(51, 4).Value = 3848
$ws.Cells.Item(52, 2).Value = 4381
$ws.Cells.Item(52, 3).Value = 515
$ws.Cells.Item(52, 4).Value = 940
$ws.Cells.Item(53, 2).Value = 15428
$ws.Cells.Item(53, 3).Value = 3194
$ws.Cells.Item(53, 4).Value = 2586
$ws.Cells.Item(54, 2).Value = 2459
$ws.Cells.Item(54, 3).Value = 791
$ws.Cells.Item(54, 4).Value = 1222
$ws.Cells.Item(55, 2).Value = 2199
$ws.Cells.Item(55, 3).Value = 1522
$ws.Cells.Item(55, 4).Value = 248
$ws.Cells.Item(56, 2).Value = 3805
$ws.Cells.Item(56, 3).Value = 1110
$ws.Cells.Item(56, 4).Value = 1276
$ws.Cells.Item(57, 2).Value = 15130
$ws.Cells.Item(57, 3).Value = 5800
$ws.Cells.Item(57, 4).Value = 3004
$ws.Cells.Item(58, 2).Value = 14413
$ws.Cells.Item(58, 3).Value = 852
$ws.Cells.Item(58, 4).Value = 501
$ws.Cells.Item(59, 2).Value = 734495
$ws.Cells.Item(59, 3).Value = 109150
$ws.Cells.Item(59, 4).Value = 103324
